$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for the affected rows
$ws.Range("F3").Value = -5
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -9
$ws.Range("F15").Value = -2
$ws.Range("F18").Value = -7
$ws.Range("F20").Value = 1
$ws.Range("F23").Value = -2
$ws.Range("F24").Value = 3
$ws.Range("F27").Value = -4
